# Add a new question row (row 14) to the "Medium" problems log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row content
$question = "Average commute time"
$difficulty = "Medium"
$url = "https://www.interviewquery.com/questions/average-commute-time"
$sqlComment = "Follows number 8 on top - for avg time across all rides in NY, use a subquery - without the groupby , a partition by can be used in select expression"

$ws.Range("C14").Value = $url
$ws.Range("A14").Value = $question
$ws.Range("B14").Value = $difficulty
$ws.Range("D14").Value = $sqlComment

# Match formatting used by the rest of the table: wrap text on C/D columns
# and a row height sized to fit the wrapped comment text.
$ws.Range("C14:D14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 68

# Update the view to reflect where the user ended up after adding the row
# (scrolled down so row 13 is at the top, with D17 selected).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select()
